$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.47685718536377
$ws.Range("B1").Value = 3.607518434524536
$ws.Range("C1").Value = 3.10010552406311
$ws.Range("D1").Value = 1.409850239753723
$ws.Range("E1").Value = 0.7790331840515137
